# Committing fixed search test cases
# Update the "Runmode" column (D) for rows 4 through 22 on the
# "Test Cases" sheet from "Y" to "N".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("D4:D22").Value = "N"
